$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 0.1429335285543671
$ws.Range("C7").Value = 1.020635186052919
$ws.Range("D7").Value = 5.36814301143687
$ws.Range("E7").Value = 2.316925335749271
$ws.Range("F7").Value = 2.345315505078564
$ws.Range("G7").Value = 36

$ws.Range("B8").Value = 0.1686795058388147
$ws.Range("C8").Value = 1.086338906570454
$ws.Range("D8").Value = 5.567529296813674
$ws.Range("E8").Value = 2.359561250913753
$ws.Range("F8").Value = 2.38788411286034
$ws.Range("G8").Value = 35

$ws.Range("B9").Value = 0.04759320168798861
$ws.Range("C9").Value = 1.466057361487795
$ws.Range("D9").Value = 8.975172192828293
$ws.Range("E9").Value = 2.995859174398605
$ws.Range("F9").Value = 3.073298772873998
$ws.Range("G9").Value = 20

$ws.Range("B10").Value = -0.6842028199372995
$ws.Range("C10").Value = 1.206520771940638
$ws.Range("D10").Value = 5.849818872833023
$ws.Range("E10").Value = 2.418639880766259
$ws.Range("F10").Value = 2.414572250977685
$ws.Range("G10").Value = 13

$ws.Range("B11").Value = 0.08624642082877883
$ws.Range("C11").Value = 0.5880280703951859
$ws.Range("D11").Value = 0.4845231414041202
$ws.Range("E11").Value = 0.6960769651440279
$ws.Range("F11").Value = 0.772240811128842
$ws.Range("G11").Value = 5
